# Insert a new weekly price record for Ciboulette (Feria Lagunitas de Puerto
# Montt) at row 115, pushing the existing rows 115-132 down to 116-133.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("115:115").Insert()

$ws.Cells.Item(115, 1).Value = 4
$ws.Cells.Item(115, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(115, 3).Value = "Los Lagos"
$ws.Cells.Item(115, 4).Value = 44504
$ws.Cells.Item(115, 5).Value = 10
$ws.Cells.Item(115, 6).Value = 100112039
$ws.Cells.Item(115, 7).Value = "Ciboulette"
$ws.Cells.Item(115, 8).Value = "Sin especificar"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 80
$ws.Cells.Item(115, 11).Value = 2500
$ws.Cells.Item(115, 12).Value = 2500
$ws.Cells.Item(115, 13).Value = 2500
$ws.Cells.Item(115, 14).Value = "`$/docena de atados"
$ws.Cells.Item(115, 15).Value = "Región Metropolitana"
$ws.Cells.Item(115, 16).Value = 833
$ws.Cells.Item(115, 17).Value = 3
$ws.Cells.Item(115, 18).Value = "Hortaliza"
